# Weekly cryptos data refresh (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "44.552.25" (thousands-dot formatted) -
# force text format so Excel does not reinterpret it as a number and
# strip formatting (e.g. "1.00" -> 1, or re-parse "2.263.68" as a date).

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '44.552.25'
$ws.Range('E2').Value = '  +1.29%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.243.46'
$ws.Range('E3').Value = '  +0.30%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +1.05%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.47'
$ws.Range('E5').Value = '  +0.02%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.18'
$ws.Range('E6').Value = '  +0.66%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  +0.60%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.12%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.523'
$ws.Range('E9').Value = '  +0.71%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.92'
$ws.Range('E10').Value = '  +0.52%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  -0.19%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.25'
$ws.Range('E12').Value = '  +0.63%  '

# Row 13
$ws.Range('E13').Value = '  +0.05%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.263.68'
$ws.Range('E14').Value = '  +1.21%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.837'
$ws.Range('E15').Value = '  +1.39%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '13.61'
$ws.Range('E16').Value = '  +0.04%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '44.269.61'
$ws.Range('E17').Value = '  +0.88%  '

# Row 18
$ws.Range('E18').Value = '  -0.64%  '

# Row 19
$ws.Range('E19').Value = '  +1.44%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.96'
$ws.Range('E20').Value = '  -1.00%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '65.60'

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '237.70'
$ws.Range('E22').Value = '  +0.65%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.96'
$ws.Range('E23').Value = '  +0.76%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.98'
$ws.Range('E24').Value = '  +1.19%  '

# Row 25
$ws.Range('E25').Value = '  +0.06%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  +3.49%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '38.12'
$ws.Range('E27').Value = '  +1.81%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.82'
$ws.Range('E28').Value = '  -1.69%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.02'
$ws.Range('E29').Value = '  +1.07%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.99'
$ws.Range('E30').Value = '  +0.71%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '153.16'
$ws.Range('E31').Value = '  +0.10%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0798'
$ws.Range('E32').Value = '  -0.28%  '

# Row 33
$ws.Range('E33').Value = '  +3.02%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.07'
$ws.Range('E34').Value = '  -5.74%  '

# Row 35
$ws.Range('E35').Value = '  +2.71%  '

# Row 36
$ws.Range('E36').Value = '  +0.73%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').Value = '  +2.59%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '14.93'
$ws.Range('E38').Value = '  -1.21%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.41'
$ws.Range('E39').Value = '  +1.83%  '

# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.77'
$ws.Range('E40').Value = '  -1.59%  '

# Row 41
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0301'
$ws.Range('E41').Value = '  +0.43%  '

# Row 42
$ws.Range('E42').Value = '  +0.19%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.793.27'
$ws.Range('E43').Value = '  +3.90%  '

# Row 44
$ws.Range('E44').Value = '  +2.72%  '

# Row 45
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.66'
$ws.Range('E45').Value = '  +11.47%  '

# Row 46
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '79.28'
$ws.Range('E46').Value = '  -7.16%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '70.63'
$ws.Range('E47').Value = '  +2.44%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '99.18'
$ws.Range('E48').Value = '  -0.74%  '

# Row 49
$ws.Range('E49').Value = '  -0.01%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.13'
$ws.Range('E50').Value = '  +0.49%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '54.66'
$ws.Range('E51').Value = '  +1.04%  '
